$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.754.25"
$ws.Range("E2").Value = "'  -0.71%  "
$ws.Range("D3").Value = "'3.092.37"
$ws.Range("E3").Value = "'  -0.80%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "'  -0.04%  "
$ws.Range("D5").Value = "'575.91"
$ws.Range("E5").Value = "'  -0.69%  "
$ws.Range("D6").Value = "'177.19"
$ws.Range("E6").Value = "'  +2.15%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("D8").Value = "'3.092.90"
$ws.Range("E8").Value = "'  -0.67%  "
$ws.Range("D9").Value = "'0.514"
$ws.Range("E9").Value = "'  -1.40%  "
$ws.Range("E10").Value = "'  -2.20%  "
$ws.Range("D11").Value = "'0.151"
$ws.Range("E11").Value = "'  -2.17%  "
$ws.Range("D12").Value = "'0.466"
$ws.Range("E12").Value = "'  -2.84%  "
$ws.Range("D13").Value = "'0.0000240"
$ws.Range("E13").Value = "'  -3.37%  "
$ws.Range("D14").Value = "'35.92"
$ws.Range("E14").Value = "'  -2.70%  "
$ws.Range("E15").Value = "'  -0.53%  "
$ws.Range("D16").Value = "'3.604.78"
$ws.Range("E16").Value = "'  -0.78%  "
$ws.Range("D17").Value = "'66.670.51"
$ws.Range("E17").Value = "'  -0.75%  "
$ws.Range("D18").Value = "'6.97"
$ws.Range("E18").Value = "'  -1.95%  "
$ws.Range("D19").Value = "'16.71"
$ws.Range("E19").Value = "'  +1.43%  "
$ws.Range("D20").Value = "'3.088.76"
$ws.Range("E20").Value = "'  -0.73%  "
$ws.Range("D21").Value = "'481.17"
$ws.Range("E21").Value = "'  -1.89%  "
$ws.Range("D22").Value = "'7.75"
$ws.Range("E22").Value = "'  -2.15%  "
$ws.Range("D23").Value = "'0.688"
$ws.Range("E23").Value = "'  -2.49%  "
$ws.Range("D24").Value = "'83.43"
$ws.Range("E24").Value = "'  -0.63%  "
$ws.Range("D25").Value = "'12.65"
$ws.Range("E25").Value = "'  -4.25%  "
$ws.Range("D26").Value = "'2.23"
$ws.Range("E26").Value = "'  -3.18%  "
$ws.Range("D27").Value = "'10.12"
$ws.Range("E27").Value = "'  -4.08%  "
$ws.Range("E28").Value = "'  +0.15%  "
$ws.Range("E29").Value = "'  -0.37%  "
$ws.Range("E30").Value = "'  -4.29%  "
$ws.Range("E31").Value = "'  -2.88%  "
$ws.Range("D32").Value = "'27.95"
$ws.Range("E32").Value = "'  -1.72%  "
$ws.Range("E33").Value = "'  -2.06%  "
$ws.Range("D34").Value = "'0.0₃0943"
$ws.Range("E34").Value = "'  -0.36%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "'  +0.05%  "
$ws.Range("D36").Value = "'48.45"
$ws.Range("E36").Value = "'  +2.83%  "
$ws.Range("D37").Value = "'5.58"
$ws.Range("E37").Value = "'  -5.04%  "
$ws.Range("D38").Value = "'0.939"
$ws.Range("E38").Value = "'  -3.58%  "
$ws.Range("B39").Value = "TheGraph"
$ws.Range("C39").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D39").Value = "'0.309"
$ws.Range("E39").Value = "'  -0.22%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").Value = "'48.98"
$ws.Range("E40").Value = "'  -2.16%  "
$ws.Range("E41").Value = "'  -2.73%  "
$ws.Range("E42").Value = "'  -0.63%  "
$ws.Range("D43").Value = "'8.31"
$ws.Range("E43").Value = "'  -2.00%  "
$ws.Range("D44").Value = "'2.68"
$ws.Range("E44").Value = "'  +2.98%  "
$ws.Range("D45").Value = "'2.786.67"
$ws.Range("E45").Value = "'  -0.80%  "
$ws.Range("D46").Value = "'370.63"
$ws.Range("E46").Value = "'  -4.14%  "
$ws.Range("D47").Value = "'135.31"
$ws.Range("E47").Value = "'  -0.01%  "
$ws.Range("E48").Value = "'  -2.71%  "
$ws.Range("D50").Value = "'24.73"
$ws.Range("E50").Value = "'  -0.99%  "
$ws.Range("D51").Value = "'2.23"
$ws.Range("E51").Value = "'  +1.02%  "
